$wb = $excel.ActiveWorkbook

# --- Notes sheet: fix "Units of measure" line ---
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A3").Value = "Units of measure: constant 2015 US$"

# --- Data sheet: fill in the data rows ---
$data = $wb.Worksheets.Item("Data")

$data.Range("A2").Value = "africa"
$data.Range("B2").Value = "Africa, regional"
$data.Range("C2").Value = 2015
$data.Range("D2").Value = 1141170000

$data.Range("A3").Value = "asia"
$data.Range("B3").Value = "Asia, regional"
$data.Range("C3").Value = 2015
$data.Range("D3").Value = 41973900

$data.Range("A4").Value = "bilateral-unspecified"
$data.Range("B4").Value = "Bilateral, unspecified"
$data.Range("C4").Value = 2015
$data.Range("D4").Value = 459683250

$data.Range("A5").Value = "middle-east"
$data.Range("B5").Value = "Middle East, regional"
$data.Range("C5").Value = 2015
$data.Range("D5").Value = 5310240000
